$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "GRT-USD"
$ws.Range("A23").Value = "BSCX-USD"
